$p = $ppt.ActivePresentation
$th = $p.SlideMaster.Theme
$tcs = $th.ThemeColorScheme
$tcs.Colors(1).RGB = 0
$tcs.Colors(2).RGB = 16777215
$tcs.Colors(3).RGB = 6968388
$tcs.Colors(4).RGB = 15132391
$tcs.Colors(5).RGB = 13998939
$tcs.Colors(6).RGB = 3243501
$tcs.Colors(7).RGB = 10855845
$tcs.Colors(8).RGB = 49407
$tcs.Colors(9).RGB = 12874308
$tcs.Colors(10).RGB = 4697456
$tcs.Colors(11).RGB = 12673797
$tcs.Colors(12).RGB = 7491477
Write-Output "done"
